$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: add C3 = "ok"
$ws.Range("C3").Value = "ok"

# Row 8: add J8 = "Out of band"
$ws.Range("J8").Value = "Out of band"

# Row 11: height changed from 57 to 38.25
$ws.Rows.Item(11).RowHeight = 38.25

# Row 18: A18 = "Chuong 1", B18 = 30
$ws.Range("A18").Value = "Chuong 1"
$ws.Range("B18").Value = 30

# Row 19: A19 = "Chuong 2", B19 = 30
$ws.Range("A19").Value = "Chuong 2"
$ws.Range("B19").Value = 30

# Row 20: A20 = "Chuong 3", B20 = 30
$ws.Range("A20").Value = "Chuong 3"
$ws.Range("B20").Value = 30

# Row 21 (new): A21 = "Phu luc", B21 = 10
$ws.Range("A21").Value = "Phu luc"
$ws.Range("B21").Value = 10

# Update selection to C21
$ws.Range("C21").Select()
